$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of bitcoin-buy data logged on 2025-08-03.
#
# Column A stores the date as a literal text string (matching every
# other "MM/DD/YYYY" row already in the sheet), not a real date
# serial. Assigning a date-shaped string straight to .Value makes
# Excel "smart" parse it into a date number/style, so instead we
# build the text via a formula (whose result is always text, never
# re-parsed) in a scratch cell, then copy only the *value* over with
# PasteSpecial - this keeps the destination cell style-free, exactly
# like the existing text-date cells.
$ws.Range("F1").Formula = "=""08/03/2025"""
$ws.Range("F1").Copy()
$ws.Range("A37").PasteSpecial(-4163) # xlPasteValues
$ws.Range("F1").ClearContents()

$ws.Range("B37").Value = 0.0004380899999999986
$ws.Range("C37").Value = 114131.7994019497
$ws.Range("D37").Value = 50
